$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: append new transmission-path rows (18-30) and expand Table1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$sheet1Rows = @(
    @("T4 child",  "T17 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T18 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T19 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T20 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T21 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T22 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T23 child", "Taringa", "Ironside State School", "Ironside State School",               "Delta (B.1.617.2)", "Wild"),
    @("T17 child", "T24",       "Taringa", "Ironside State School", "Ironside State School Close Contact", "Delta (B.1.617.2)", "Wild"),
    @("T18 child", "T25",       "Taringa", "Ironside State School", "Ironside State School Close Contact", "Delta (B.1.617.2)", "Wild"),
    @("T19 child", "T26",       "Taringa", "Ironside State School", "Ironside State School Close Contact", "Delta (B.1.617.2)", "Wild"),
    @("T20 child", "T27",       "Taringa", "Ironside State School", "Ironside State School Close Contact", "Delta (B.1.617.2)", "Wild"),
    @("T21 child", "T28",       "Taringa", "Ironside State School", "Ironside State School Close Contact", "Delta (B.1.617.2)", "Wild"),
    @("T4 child",  "T29",       "Taringa", "Karate Class",          "Karate Class",                        "Delta (B.1.617.2)", "Wild")
)

$r = 18
foreach ($row in $sheet1Rows) {
    # Copy the date cell immediately above so the new row picks up the same
    # date-number format (style index) instead of creating a brand-new style.
    $ws1.Range("A" + ($r - 1)).Copy($ws1.Range("A" + $r))
    $ws1.Range("A" + $r).Value = 44410

    $ws1.Range("B" + $r).Value = $row[0]
    $ws1.Range("C" + $r).Value = $row[1]
    $ws1.Range("D" + $r).Value = $row[2]
    $ws1.Range("E" + $r).Value = $row[3]
    $ws1.Range("F" + $r).Value = $row[4]
    $ws1.Range("G" + $r).Value = $row[5]
    $ws1.Range("H" + $r).Value = $row[6]

    $r = $r + 1
}

# Expand Table1 (and its AutoFilter) to cover the newly added rows.
$table1 = $ws1.ListObjects.Item(1)
$table1.Resize($ws1.Range("A1:H30"))

$ws1.Range("H30").Select()

# ---------------------------------------------------------------------------
# "Date Colours" sheet: refresh the colour-gradient lookup values
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Date Colours")

$ws2.Range("B2").Value = "#f0e2f0"
$ws2.Range("B3").Value = "#e0c6e1"
$ws2.Range("B4").Value = "#d0aad2"
$ws2.Range("E4").Value = "#e0c1c1"
$ws2.Range("B5").Value = "#c08ec3"

$ws2.Range("E4:E5").Select()

# Sheet1 becomes the active/visible tab again (it was "Date Colours" before).
$ws1.Activate()
